$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "Förändrad" (Changed) date column (C) for rows 2-23
# from serial 45182 (2023-09-13) to 45184 (2023-09-15)
for ($row = 2; $row -le 23; $row++) {
    $cell = $ws.Cells.Item($row, 3)
    if ($cell.Value2 -eq 45182) {
        $cell.Value2 = 45184
    }
}
